$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing (default) style of the data cells in column D so we can
# restore it after using a quote-prefix to force text entry (prevents Excel from
# auto-converting numeric-looking strings like "60.23" into floating point numbers).
$styleD = $ws.Range("D2").Style

$ws.Range("D2").Value = "'38.740.90"
$ws.Range("E2").Value = "  +2.44%  "
$ws.Range("D3").Value = "'2.085.69"
$ws.Range("E3").Value = "  +1.92%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'228.32"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("E6").Value = "  +0.69%  "
$ws.Range("D7").Value = "'60.23"
$ws.Range("E7").Value = "  +1.08%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +1.98%  "
$ws.Range("E10").Value = "  +0.55%  "
$ws.Range("E11").Value = "  -0.88%  "
$ws.Range("D12").Value = "'2.396.18"
$ws.Range("E12").Value = "  +2.01%  "
$ws.Range("D13").Value = "'14.97"
$ws.Range("E13").Value = "  +3.93%  "
$ws.Range("D14").Value = "'21.85"
$ws.Range("E14").Value = "  +2.02%  "
$ws.Range("D15").Value = "'0.794"
$ws.Range("E15").Value = "  +4.08%  "
$ws.Range("D16").Value = "'5.49"
$ws.Range("E16").Value = "  -0.32%  "
$ws.Range("D17").Value = "'2.086.95"
$ws.Range("E17").Value = "  +2.31%  "
$ws.Range("D18").Value = "'38.683.82"
$ws.Range("E18").Value = "  +2.40%  "
$ws.Range("D19").Value = "'71.45"
$ws.Range("E19").Value = "  +2.81%  "
$ws.Range("D20").Value = "'6.03"
$ws.Range("E20").Value = "  +1.90%  "
$ws.Range("E21").Value = "  +1.15%  "
$ws.Range("D22").Value = "'227.09"
$ws.Range("E22").Value = "  +2.06%  "
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("D24").Value = "'2.38"
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("D25").Value = "'2.33"
$ws.Range("E25").Value = "  +2.32%  "
$ws.Range("D26").Value = "'170.70"
$ws.Range("E26").Value = "  +1.06%  "
$ws.Range("D27").Value = "'9.51"
$ws.Range("E27").Value = "  +2.40%  "
$ws.Range("D28").Value = "'0.138"
$ws.Range("E28").Value = "  +8.44%  "
$ws.Range("E29").Value = "  +13.02%  "
$ws.Range("D30").Value = "'19.15"
$ws.Range("E30").Value = "  +1.92%  "
$ws.Range("E31").Value = "  +0.97%  "
$ws.Range("E32").Value = "  +5.96%  "
$ws.Range("D33").Value = "'4.50"
$ws.Range("E33").Value = "  +2.85%  "
$ws.Range("D34").Value = "'4.68"
$ws.Range("E34").Value = "  +3.43%  "
$ws.Range("E35").Value = "  +0.93%  "
$ws.Range("D36").Value = "'6.46"
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("E37").Value = "  +1.30%  "
$ws.Range("E38").Value = "  +2.23%  "
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("D40").Value = "'17.89"
$ws.Range("E40").Value = "  -2.65%  "
$ws.Range("E41").Value = "  +4.91%  "
$ws.Range("D42").Value = "'1.539.78"
$ws.Range("E42").Value = "  +1.05%  "
$ws.Range("D43").Value = "'100.60"
$ws.Range("E43").Value = "  +3.16%  "
$ws.Range("E44").Value = "  -0.74%  "
$ws.Range("E45").Value = "  +3.45%  "
$ws.Range("D46").Value = "'7.68"
$ws.Range("E46").Value = "  +8.16%  "
$ws.Range("E47").Value = "  +1.30%  "
$ws.Range("D48").Value = "'4.11"
$ws.Range("E48").Value = "  -1.66%  "
$ws.Range("E49").Value = "  +2.48%  "
$ws.Range("E50").Value = "  +0.65%  "
$ws.Range("D51").Value = "'2.284.58"

# Restore the original (default/unstyled) formatting to the column-D cells we
# touched, so only the cell *values* change -- not their style/format.
$ws.Range("D2:D51").Style = $styleD

